$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new student row (row 4) following the same pattern as rows 2-3
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 200155
$ws.Range("C4").Value = "Mohamed Elhefny"
$ws.Range("D4:O4").Value = " "
